$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.153.85'
$ws.Range("E2").Value = '  +2.54%  '
$ws.Range("D3").Value = '3.062.99'
$ws.Range("E3").Value = '  +2.65%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '526.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.84%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.449'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.59'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.28%  '
$ws.Range("E10").Value = '  +7.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.371'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.06%  '
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").Value = '3.594.02'
$ws.Range("E13").Value = '  +2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.40'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.06%  '
$ws.Range("E15").Value = '  +15.74%  '
$ws.Range("D16").Value = '58.150.10'
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("E17").Value = '  +10.02%  '
$ws.Range("D18").Value = '3.072.11'
$ws.Range("E18").Value = '  +2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '342.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.96%  '
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.74'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  +8.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.56%  '
$ws.Range("D26").Value = '0.0₃0972'
$ws.Range("E26").Value = '  +9.77%  '
$ws.Range("E27").Value = '  +3.62%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.81%  '
$ws.Range("E30").Value = '  +10.02%  '
$ws.Range("E31").Value = '  +7.28%  '
$ws.Range("E32").Value = '  +6.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.36'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.97'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.71%  '
$ws.Range("E37").Value = '  +3.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.68%  '
$ws.Range("E39").Value = '  +4.37%  '
$ws.Range("D40").Value = '3.100.70'
$ws.Range("E40").Value = '  +2.61%  '
$ws.Range("E41").Value = '  +3.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +11.89%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.667'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.27%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.09%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.340.55'
$ws.Range("E46").Value = '  +4.95%  '
$ws.Range("E47").Value = '  +4.64%  '
$ws.Range("E48").Value = '  +4.26%  '
$ws.Range("E49").Value = '  +6.31%  '
$ws.Range("E50").Value = '  +3.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.93%  '
